# Applies the "packets_counter" workbook update:
#  - Sheet "R3" is renamed to "R2" (data refreshed: em0/em1/em2 rx/tx octet & unicast counters)
#  - Sheet "SW1" is renamed to "R3" (interfaces renumbered Ethernet1/x, FastEthernet0/x rows,
#    counters refreshed, trailing rows trimmed to 12)
#  - Sheet "R1" keeps its name, counters refreshed, trailing Loopback0 row removed (12 rows)
#  - A brand-new sheet "MLS1" is appended with its own interface/counter table
#  - The newly appended sheet becomes the active tab

function Set-RowValues {
    param($ws, $rowNum, $values)
    $n = $values.Count
    $arr = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0, $i] = $values[$i]
    }
    $lastCol = [char](64 + $n)
    $addr = "A${rowNum}:${lastCol}${rowNum}"
    $ws.Range($addr).Value = $arr
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "R1" - update counters, drop the Loopback0 row (13 -> 12 rows)
# ---------------------------------------------------------------------------
$wsR1 = $wb.Worksheets.Item(1)

Set-RowValues $wsR1 3  @("Ethernet1/0",     0,   0, 0, 0, 1435055, 15710, 0, 0, 875237, 7636)
Set-RowValues $wsR1 11 @("FastEthernet0/0", 337, 0, 0, 0, 33276,   352,   0, 0, 322402, 2941)
Set-RowValues $wsR1 12 @("FastEthernet0/1", 3,   0, 0, 0, 884,     10,    0, 0, 280803, 2501)
$wsR1.Rows.Item(13).Delete()

# ---------------------------------------------------------------------------
# Sheet 2: old "R3" -> renamed to "R2" - refresh em0/em1/em2 counters
# ---------------------------------------------------------------------------
$wsR2 = $wb.Worksheets.Item(2)
$wsR2.Name = "R2"

Set-RowValues $wsR2 3 @("em0", -1, 0, 0, -1, 514565, -1, 0, 0, 922454, -1)
Set-RowValues $wsR2 4 @("em1", -1, 0, 0, -1, 31734,  -1, 0, 0, 33276,  -1)
Set-RowValues $wsR2 5 @("em2", -1, 0, 0, -1, 31710,  -1, 0, 0, 32912,  -1)
# row 6 ("mtun") is unchanged

# ---------------------------------------------------------------------------
# Sheet 3: old "SW1" -> renamed to "R3" - renumber interfaces, refresh counters,
# drop the legacy GigabitEthernet2/1..Vlan10 tail and replace row 12 with the new
# FastEthernet0/1 entry (19 -> 12 rows)
# ---------------------------------------------------------------------------
$wsR3 = $wb.Worksheets.Item(3)
$wsR3.Name = "R3"

Set-RowValues $wsR3 3  @("Ethernet1/0",     0,   0, 0, 0, 1533184, 16734, 0, 0, 1268448, 11013)
Set-RowValues $wsR3 4  @("Ethernet1/1",     0,   0, 0, 0, 0,       0,     0, 0, 0,       0)
Set-RowValues $wsR3 5  @("Ethernet1/2",     0,   0, 0, 0, 0,       0,     0, 0, 0,       0)
Set-RowValues $wsR3 6  @("Ethernet1/3",     0,   0, 0, 0, 0,       0,     0, 0, 0,       0)
Set-RowValues $wsR3 7  @("Ethernet1/4",     0,   0, 0, 0, 0,       0,     0, 0, 0,       0)
Set-RowValues $wsR3 8  @("Ethernet1/5",     0,   0, 0, 0, 0,       0,     0, 0, 0,       0)
Set-RowValues $wsR3 9  @("Ethernet1/6",     0,   0, 0, 0, 0,       0,     0, 0, 0,       0)
Set-RowValues $wsR3 10 @("Ethernet1/7",     0,   0, 0, 0, 0,       0,     0, 0, 0,       0)
Set-RowValues $wsR3 11 @("FastEthernet0/0", 333, 0, 0, 0, 32912,   348,   0, 0, 324353,  2966)

# Drop the old rows 13-19 (GigabitEthernet2/1 .. Vlan10), keep row 12 for now
$wsR3.Range("A13:K19").EntireRow.Delete()

# Row 12 becomes the new FastEthernet0/1 entry
Set-RowValues $wsR3 12 @("FastEthernet0/1", 0, 0, 0, 0, 0, 0, 0, 0, 320372, 2978)

# ---------------------------------------------------------------------------
# Sheet 4 (new): "MLS1"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMLS1 = $wb.Worksheets.Add($null, $lastSheet)
$wsMLS1.Name = "MLS1"

Set-RowValues $wsMLS1 2 @("interface", "rx_broadcast", "rx_discards", "rx_errors", "rx_multicast", "rx_octets", "rx_unicast", "tx_discards", "tx_errors", "tx_octets", "tx_unicast")

Set-RowValues $wsMLS1 3  @("GigabitEthernet0/0", 745, 0, 0, 0, 982276, 10002, 0, 0, 1670165, 13530)
Set-RowValues $wsMLS1 4  @("GigabitEthernet0/1", 4,   0, 0, 0, 0,      0,     0, 0, 62177,   632)
Set-RowValues $wsMLS1 5  @("GigabitEthernet0/2", 3,   0, 0, 3, 226,    3,     0, 0, 90898,   1106)
Set-RowValues $wsMLS1 6  @("GigabitEthernet0/3", 3,   0, 0, 3, 226,    3,     0, 0, 91062,   1108)
Set-RowValues $wsMLS1 7  @("GigabitEthernet1/0", 0,   0, 0, 0, 0,      0,     0, 0, 0,       0)
Set-RowValues $wsMLS1 8  @("GigabitEthernet1/1", 0,   0, 0, 0, 0,      0,     0, 0, 0,       0)
Set-RowValues $wsMLS1 9  @("GigabitEthernet1/2", 0,   0, 0, 0, 0,      0,     0, 0, 0,       0)
Set-RowValues $wsMLS1 10 @("GigabitEthernet1/3", 0,   0, 0, 0, 0,      0,     0, 0, 0,       0)
Set-RowValues $wsMLS1 11 @("GigabitEthernet2/0", 0,   0, 0, 0, 0,      0,     0, 0, 0,       0)
Set-RowValues $wsMLS1 12 @("GigabitEthernet2/1", 0,   0, 0, 0, 0,      0,     0, 0, 0,       0)
Set-RowValues $wsMLS1 13 @("GigabitEthernet2/2", 0,   0, 0, 0, 0,      0,     0, 0, 0,       0)
Set-RowValues $wsMLS1 14 @("GigabitEthernet2/3", 0,   0, 0, 0, 0,      0,     0, 0, 0,       0)
Set-RowValues $wsMLS1 15 @("GigabitEthernet3/0", 0,   0, 0, 0, 0,      0,     0, 0, 0,       0)
Set-RowValues $wsMLS1 16 @("GigabitEthernet3/1", 0,   0, 0, 0, 0,      0,     0, 0, 0,       0)
Set-RowValues $wsMLS1 17 @("GigabitEthernet3/2", 0,   0, 0, 0, 0,      0,     0, 0, 0,       0)
Set-RowValues $wsMLS1 18 @("Vlan30",             0,   0, 0, 0, 86,     1,     0, 0, 53032,   634)
Set-RowValues $wsMLS1 19 @("Vlan40",             0,   0, 0, 0, 86,     1,     0, 0, 53196,   636)
Set-RowValues $wsMLS1 20 @("Vlan50",             0,   0, 0, 0, 0,      0,     0, 0, 53048,   634)

# Newly added sheet is the active tab (matches activeTab="3")
$wsMLS1.Activate()
